$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# Update the confidential disclosure date from 2021-03-22 to 2021-03-23
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-23 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.1101264281381614
$ws.Range("E2").Value = -0.03874300473525616

$ws.Range("D3").Value = 0.1047021675560621
$ws.Range("E3").Value = -0.005222906174221276

$ws.Range("D4").Value = 0.1149195414895816
$ws.Range("E4").Value = -0.006706408345752646

$ws.Range("D5").Value = 0.1363478329474702
$ws.Range("E5").Value = -0.01573426573426573

$ws.Range("D6").Value = 0.1312848002434956
$ws.Range("E6").Value = -0.001023541453428978

$ws.Range("D7").Value = 0.1455741671802382
$ws.Range("E7").Value = -0.02456174147563073

$ws.Range("D8").Value = 0.1280405933066407
$ws.Range("E8").Value = -0.02428785607196404

$ws.Range("D9").Value = 0.1290044691383502
$ws.Range("E9").Value = -0.01686030124286853

$ws.Range("E10").Value = -0.01672432494025411
